$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $range as a genuine text (shared-string) cell,
# no matter whether it "looks like" a number. Routing the literal through
# a =TEXT(...,"@") formula and then pasting-as-values forces Excel to
# store the result as a string (t="s") instead of re-parsing it as a
# number - and unlike the "leading apostrophe" trick it does NOT leave a
# quotePrefix cell style behind in styles.xml.
function Set-TextValue($range, $text) {
    $range.Formula = '=TEXT("' + $text + '","@")'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# --- 1. Capture the existing data-row values before they are moved -------
$oldA2 = $ws.Range("A2").Text
$oldB2 = $ws.Range("B2").Text
$oldC2 = $ws.Range("C2").Text
$oldD2 = $ws.Range("D2").Text

# --- 2. Clear the cells/rows being restructured ---------------------------
$ws.Range("A2").Clear()
$ws.Range("B2").Clear()
$ws.Range("C2").Clear()
$ws.Range("D2").Clear()
$ws.Range("A3").Clear()

# --- 3. Rewrite the header row (row 1) -------------------------------------
$ws.Range("A1").Value = "Unnamed: 0"
Set-TextValue $ws.Range("B1") "2019"
$ws.Range("C1").Value = "Unnamed: 1"

# Give the new header cells (B1, C1) the same look as A1 (bold/bordered).
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Write the reshaped data into column B, rows 2-5 --------------------
Set-TextValue $ws.Range("B2") $oldA2
Set-TextValue $ws.Range("B3") $oldB2
Set-TextValue $ws.Range("B4") $oldC2
Set-TextValue $ws.Range("B5") $oldD2
$excel.CutCopyMode = $false
